# Cloudflare DNS records export was refreshed; 11 stale DNS record rows
# (previously rows 105-115, covering old italiadns.com mailtrap CNAME/MX
# records that no longer exist) were removed from the export. Deleting
# these rows shifts every row below them up by 11, so the former row 116
# becomes row 105, and the sheet shrinks from 156 data rows to 145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 105 through 115 (inclusive) - this shifts all subsequent
# rows up by 11 automatically, matching the target diff exactly.
$ws.Range("A105:O115").EntireRow.Delete()
